$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header style (bold, centered, bordered) by copying
# the format from the neighboring header cell (H1) instead of creating a
# brand new style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..29 for columns I (I0) and J (IF)
$data = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 7, 7),
    @(5, 6, 6),
    @(6, 4, 6),
    @(7, 6, 6),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 8, 8),
    @(11, 8, 8),
    @(12, 6, 6),
    @(13, 3, 3),
    @(14, 7, 8),
    @(15, 9, 9),
    @(16, 6, 7),
    @(17, 8, 8),
    @(18, 7, 8),
    @(19, 9, 9),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 7, 7),
    @(23, 5, 5),
    @(24, 4, 4),
    @(25, 6, 6),
    @(26, 8, 8),
    @(27, 7, 7),
    @(28, 7, 7),
    @(29, 4, 4)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
